# Update the Clinic Template: add a header row describing the columns
# used by the clinic staffing sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header values (row 1) -------------------------------------------------
$ws.Range("A1").Value = "Clinic Name"
$ws.Range("B1").Value = "Min Number of Pediatric Doctors Needed "
$ws.Range("C1").Value = "Min Number of Adult Doctors Needed"
$ws.Range("D1").Value = "Ideal # of Providers"
$ws.Range("E1").Value = "Max # of Providers"

# --- Formatting --------------------------------------------------------
# Whole header row uses 10pt Arial instead of the workbook default font.
$ws.Range("A1:E1").Font.Name = "Arial"
$ws.Range("A1:E1").Font.Size = 10

# The descriptive column headers (B:E) wrap their text so the long labels
# are readable; the clinic-name header (A) stays on a single line.
$ws.Range("B1:E1").WrapText = $true

# Row is tall enough to show the wrapped header text.
$ws.Rows.Item(1).RowHeight = 57

# Leave the selection where the original author left it when they saved.
$ws.Range("G2").Select() | Out-Null
